# edit.ps1 - applies the Day_14.pptx commit:
#   "add rdbms new data from lab/drive"
#
# Concrete changes reproduced here:
#   1. Delete the trailing "PL / SQL" slide (slide 5), which removes the
#      corresponding <p:sldId> entry from the presentation's slide list.
#   2. Resize the second table on the "Example - 3" slide (slide 4) so its
#      frame height grows from 436309 EMU to 456565 EMU.

$p = $ppt.ActivePresentation

# --- 1. Delete the last slide ("PL / SQL") -------------------------------
$lastIndex = $p.Slides.Count
$p.Slides.Item($lastIndex).Delete()

# --- 2. Grow the table frame on slide 4 ("Example - 3") ------------------
$slide4 = $p.Slides.Item(4)
$table3 = $slide4.Shapes.Item("Table 3")

# Target height is 456565 EMU; PowerPoint's Shape.Height is in points
# (1 pt = 12700 EMU), so convert precisely.
$targetEmu = 456565
$table3.Height = $targetEmu / 914400 * 72

$p.Save()
